$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.80898833333333
$ws.Range("H2").Value = 35.426965
$ws.Range("I2").Value = 0.08059095716837197
$ws.Range("J2").Value = 0.08059095716837197
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 910.3492816655611
$ws.Range("R2").Value = 8193.143534990051
$ws.Range("S2").Value = 0.01937258865929369
$ws.Range("T2").Value = 0.01937258865929369
$ws.Range("G3").Value = 11.80898833333333
$ws.Range("H3").Value = 35.426965
$ws.Range("I3").Value = 0.08059095716837197
$ws.Range("J3").Value = 0.08059095716837197
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 1199.557475768898
$ws.Range("R3").Value = 10796.01728192008
$ws.Range("S3").Value = 0.02552705210986123
$ws.Range("T3").Value = 0.02552705210986122
$ws.Range("G4").Value = 11.80898833333333
$ws.Range("H4").Value = 35.426965
$ws.Range("I4").Value = 0.08059095716837197
$ws.Range("J4").Value = 0.08059095716837197
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 1677.192698258125
$ws.Range("R4").Value = 15094.73428432312
$ws.Range("S4").Value = 0.03569131639921707
$ws.Range("T4").Value = 0.03569131639921707
$ws.Range("I5").Value = 0.8141849724511824
$ws.Range("J5").Value = 0.8141849724511824
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 9196.971110112472
$ws.Range("R5").Value = 82772.73999101225
$ws.Range("S5").Value = 0.1957151412275967
$ws.Range("T5").Value = 0.1957151412275967
$ws.Range("I6").Value = 0.8141849724511824
$ws.Range("J6").Value = 0.8141849724511824
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("Q6").Value = 12118.75010147916
$ws.Range("S6").Value = 0.2578917405758752
$ws.Range("T6").Value = 0.2578917405758752
$ws.Range("I7").Value = 0.8141849724511824
$ws.Range("J7").Value = 0.8141849724511824
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.3605780906477105
$ws.Range("T7").Value = 0.3605780906477105
$ws.Range("I8").Value = 0.1052240703804457
$ws.Range("J8").Value = 0.1052240703804457
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 1188.603042455963
$ws.Range("R8").Value = 10697.42738210367
$ws.Range("S8").Value = 0.02529393748579201
$ws.Range("T8").Value = 0.02529393748579201
$ws.Range("I9").Value = 0.1052240703804457
$ws.Range("J9").Value = 0.1052240703804457
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.03332954989232327
$ws.Range("T9").Value = 0.03332954989232326
$ws.Range("I10").Value = 0.1052240703804457
$ws.Range("J10").Value = 0.1052240703804457
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.04660058300233039
$ws.Range("T10").Value = 0.04660058300233038
